$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix B34: was stored as text "1", should be numeric 1 ---
$ws.Range("B34").Value = 1

# --- Add new row 35 ---
$ws.Range("A35").Value = "Ying Tang"

# B35 keeps the "typed as text" quirk seen elsewhere in the sheet: format the
# cell as Text first so the digit string isn't coerced into a number.
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "3"

$ws.Range("C35").Value = "无"
$ws.Range("D35").Value = "ACK"
$ws.Range("E35").Value = "OTH"
$ws.Range("F35").Value = "d4ad31e6-de82-4ee8-af90-c18d97ed2c36"
$ws.Range("G35").Value = "Bk7wvW-C-_annotated.xlsx"
$ws.Range("H35").Value = "We will update our paper very soon."
